# Applies the edit described by the diff:
#  - Column E (rows 17-51) holds the "Periodo Mora" values for worker
#    85164695 / SAMUEL ENRIQUE ROCHA RIVERA. Their order is reversed
#    (previously ascending 2107..2405, now descending 2405..2107).
#  - As a side effect of the underlying data reshuffle, the "Valor Mora"
#    values in column F for the first (row 17) and last (row 51) periods
#    are swapped (26650 <-> 36341); all the other F values in between are
#    identical (36341) so no visible change happens for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @(
    "2405", "2404", "2403", "2402", "2401",
    "2312", "2311", "2310", "2309", "2308", "2307", "2306", "2305", "2304", "2303", "2302", "2301",
    "2212", "2211", "2210", "2209", "2208", "2207", "2206", "2205", "2204", "2203", "2202", "2201",
    "2112", "2111", "2110", "2109", "2108", "2107"
)

$startRow = 17
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Swap the "Valor Mora" figures for the first and last period rows.
$ws.Range("F17").Value = 26650
$ws.Range("F51").Value = 36341
